$d = $word.ActiveDocument

function FindSet($doc, $old, $new) {
    $rng = $doc.Content
    $found = $rng.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Output "NOT FOUND: $old"
        return
    }
    $rng.Text = $new
}

# Title
FindSet $d "Unveiling the Enigma of Black Holes" "Exploring Biology: Unveiling the Secrets of Life"

# Author
FindSet $d "Richard Dawkins" "Harmony Kyle"

# Email local-part and domain
FindSet $d "rdawkins@evolutionarybiology" "hkyle@biosphere"
FindSet $d "org" "academy"

# Paragraph 1 (intro) sentences
FindSet $d "In the vast expanse of the cosmos, there exist celestial enigmas that capture the imagination and boggle the mind" "Biology, the science of life, embarks on an extraordinary journey through the captivating realm of living organisms"

FindSet $d " Among these celestial wonders, black holes stand as beacons of extreme gravity, marking the boundaries of our understanding of physics" " From the tiniest bacteria to the majestic blue whales, biology delves into the intricate tapestry of life's processes, revealing the underlying mechanisms that govern our existence"

# Merge of 3 runs into 1
FindSet $d " These cosmic behemoths, formed by the collapse of massive stars or the remnants of supernovae, hold secrets that have captivated scientists, philosophers, and science fiction enthusiasts alike. From their elusive nature to their potential role in shaping the universe, black holes have inspired both awe and speculation, propelling humankind's fascination with the cosmos" " As we unravel the enigmatic secrets hidden within the building blocks of life, we gain a profound appreciation for the interconnectedness of all living things and the awe-inspiring symphony of life's symphony"

FindSet $d "Journey into the depths of a black hole, and you'll find a singularity--a point of infinite density where space and time become distorted beyond recognition" "Venturing into the microscopic realm, biology unveils a captivating world of cells, the fundamental units of life"

FindSet $d " Beyond the event horizon, the boundary beyond which nothing, not even light, can escape, lies a realm of mystery and paradox" " Within these microscopic chambers, intricate processes unfold,orchestrated by complex molecules, each playing a vital role in the symphony of life"

# Merge of 3 runs into 1
FindSet $d " Here, the laws of physics, as we know them, break down, leaving us with more questions than answers. Theoreticians have proposed theories and hypotheses to explain the phenomena associated with black holes, including the concept of Hawking radiation, which suggests that black holes emit particles over time, leading to their eventual evaporation" " Discover the marvels of photosynthesis, where plants harness the energy of the sun to create sustenance, or witness the fascinating process of cell division, where life perpetuates itself"

FindSet $d "But the enigma of black holes extends beyond their theoretical implications" "Biology extends beyond the confines of individual organisms, delving into the captivating world of ecosystems, where intricate webs of interactions weave together diverse species"

FindSet $d " Their existence challenges our understanding of reality and forces us to confront fundamental questions about the nature of space, time, and gravity" " Explore the delicate balance of predator-prey relationships, or unravel the intricate network of symbiotic partnerships that sustain life on Earth"

FindSet $d " Are black holes truly portals to other dimensions, as some physicists speculate? Do they hold the key to understanding the origin and fate of the universe? As we delve deeper into the mysteries surrounding black holes, we embark on a voyage of scientific exploration that pushes the boundaries of human knowledge and understanding" " Biology unveils the interdependence of all living things, highlighting the fragility and resilience of the ecosystems we inhabit"

# Summary paragraph
FindSet $d "Black holes, enigmatic cosmic entities born from the collapse of massive stars, captivate scientists and ignite awe among humankind" "Biology, the science of life, unveils the intricacies of the living world, taking us on an awe-inspiring journey through the tapestry of life's processes"

FindSet $d " Beyond their event horizon lies a realm of mystery and paradox, where space and time distort, and the laws of physics falter" " From the symphony of cells to the interconnectedness of ecosystems, biology reveals the profound beauty and interconnectedness of all living things"

# Merge of 4 runs (incl. lastRenderedPageBreak run + trailing period run) into 1
FindSet $d " The study of black holes offers a tantalizing glimpse into the deepest secrets of the universe, propelling us toward a greater comprehension of reality and our place within it. While black holes continue to puzzle and enthrall us, they serve as constant reminders of the infinite mysteries that await discovery in the vast cosmic tapestry" " Through the exploration of biology, we gain an appreciation for the marvel of life and the importance of preserving the delicate balance of our planet's ecosystems"

# Add a new empty paragraph after the Summary paragraph (before sectPr)
$last = $d.Paragraphs.Last
$last.Range.InsertParagraphAfter()
